$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad" / last-changed date) for every existing
#    data row (2-534) from 2023-09-19 (45188) to 2023-09-20 (45189).
$ws.Range("C2:C534").Value = 45189

# 2. Row 534 was missing the standard 15pt custom row height that every
#    other data row already has; set it explicitly so it matches.
$ws.Rows.Item(534).RowHeight = 15

# 3. Append the new record as row 535.
$ws.Cells.Item(535, 1).Value = "A 44298-2023"

$ws.Cells.Item(535, 2).Value = 45187
$ws.Cells.Item(535, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(535, 3).Value = 45189
$ws.Cells.Item(535, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(535, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(535, 5).Value = "MOTALA"

$ws.Cells.Item(535, 7).Value = 2.5
$ws.Cells.Item(535, 8).Value = 0
$ws.Cells.Item(535, 9).Value = 0
$ws.Cells.Item(535, 10).Value = 0
$ws.Cells.Item(535, 11).Value = 0
$ws.Cells.Item(535, 12).Value = 0
$ws.Cells.Item(535, 13).Value = 0
$ws.Cells.Item(535, 14).Value = 0
$ws.Cells.Item(535, 15).Value = 0
$ws.Cells.Item(535, 16).Value = 0
$ws.Cells.Item(535, 17).Value = 0

# Column R keeps the same wrap-text formatting used throughout the column,
# even though this new row has no species names listed yet.
$ws.Cells.Item(535, 18).WrapText = $true
